$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($w in 46.0, 46.5, 46.875, 46.86, 46.88, 46.9, 47.0) {
  $ws.Columns("C").ColumnWidth = $w
  Write-Output ("set=" + $w + " got=" + $ws.Columns("C").ColumnWidth)
}
